$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 138, shifting existing rows 138-222 down to 139-223
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with this week's new record
$ws.Range("A138").Value = 9
$ws.Range("B138").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C138").Value = "Metropolitana"
$ws.Range("D138").Value = 44567
$ws.Range("E138").Value = 13
$ws.Range("F138").Value = 100112001
$ws.Range("G138").Value = "Berenjena"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 97
$ws.Range("K138").Value = 7000
$ws.Range("L138").Value = 8000
$ws.Range("M138").Value = 7495
$ws.Range("N138").Value = "$/caja 50 unidades"
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("P138").Value = 150
$ws.Range("Q138").Value = 50
$ws.Range("R138").Value = "Hortaliza"
